$d = $word.ActiveDocument

# Update the date title
$d.Content.Find.Execute("2025-07-17 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-18 Friday", 2) | Out-Null

# Update the math problems table
$tbl = $d.Tables.Item(1)

$tbl.Rows.Item(1).Cells.Item(1).Range.Text = "50+26="
$tbl.Rows.Item(1).Cells.Item(2).Range.Text = "73+16="
$tbl.Rows.Item(1).Cells.Item(3).Range.Text = "33-4="
$tbl.Rows.Item(1).Cells.Item(4).Range.Text = "53+34="
$tbl.Rows.Item(1).Cells.Item(5).Range.Text = "2+36="

$tbl.Rows.Item(2).Cells.Item(1).Range.Text = "42+37="
$tbl.Rows.Item(2).Cells.Item(2).Range.Text = "56-22="
$tbl.Rows.Item(2).Cells.Item(3).Range.Text = "98-48="
$tbl.Rows.Item(2).Cells.Item(4).Range.Text = "71+16="
$tbl.Rows.Item(2).Cells.Item(5).Range.Text = "2+87="

$tbl.Rows.Item(3).Cells.Item(1).Range.Text = "80-29="
$tbl.Rows.Item(3).Cells.Item(2).Range.Text = "47-42="
$tbl.Rows.Item(3).Cells.Item(3).Range.Text = "18+15="
$tbl.Rows.Item(3).Cells.Item(4).Range.Text = "4+34="
$tbl.Rows.Item(3).Cells.Item(5).Range.Text = "66-52="

$tbl.Rows.Item(4).Cells.Item(1).Range.Text = "7+47="
$tbl.Rows.Item(4).Cells.Item(2).Range.Text = "79-61="
$tbl.Rows.Item(4).Cells.Item(3).Range.Text = "98-20="
$tbl.Rows.Item(4).Cells.Item(4).Range.Text = "20-6="
$tbl.Rows.Item(4).Cells.Item(5).Range.Text = "93-87="

$tbl.Rows.Item(5).Cells.Item(1).Range.Text = "96-45="
$tbl.Rows.Item(5).Cells.Item(2).Range.Text = "31+17="
$tbl.Rows.Item(5).Cells.Item(3).Range.Text = "47-31="
$tbl.Rows.Item(5).Cells.Item(4).Range.Text = "91-19="
$tbl.Rows.Item(5).Cells.Item(5).Range.Text = "11+80="

$tbl.Rows.Item(6).Cells.Item(1).Range.Text = "61+12="
$tbl.Rows.Item(6).Cells.Item(2).Range.Text = "22-12="
$tbl.Rows.Item(6).Cells.Item(3).Range.Text = "41+58="
$tbl.Rows.Item(6).Cells.Item(4).Range.Text = "22+35="
$tbl.Rows.Item(6).Cells.Item(5).Range.Text = "64+8="

$tbl.Rows.Item(7).Cells.Item(1).Range.Text = "71+17="
$tbl.Rows.Item(7).Cells.Item(2).Range.Text = "53+39="
$tbl.Rows.Item(7).Cells.Item(3).Range.Text = "34+24="
$tbl.Rows.Item(7).Cells.Item(4).Range.Text = "3+32="
$tbl.Rows.Item(7).Cells.Item(5).Range.Text = "79-72="

$tbl.Rows.Item(8).Cells.Item(1).Range.Text = "34+64="
$tbl.Rows.Item(8).Cells.Item(2).Range.Text = "68-15="
$tbl.Rows.Item(8).Cells.Item(3).Range.Text = "29-1="
$tbl.Rows.Item(8).Cells.Item(4).Range.Text = "45-16="
$tbl.Rows.Item(8).Cells.Item(5).Range.Text = "61+5="

$tbl.Rows.Item(9).Cells.Item(1).Range.Text = "75-59="
$tbl.Rows.Item(9).Cells.Item(2).Range.Text = "59-0="
$tbl.Rows.Item(9).Cells.Item(3).Range.Text = "38+25="
$tbl.Rows.Item(9).Cells.Item(4).Range.Text = "63-37="
$tbl.Rows.Item(9).Cells.Item(5).Range.Text = "99-13="

$tbl.Rows.Item(10).Cells.Item(1).Range.Text = "19+62="
$tbl.Rows.Item(10).Cells.Item(2).Range.Text = "84-66="
$tbl.Rows.Item(10).Cells.Item(3).Range.Text = "33+7="
$tbl.Rows.Item(10).Cells.Item(4).Range.Text = "8+65="
$tbl.Rows.Item(10).Cells.Item(5).Range.Text = "87-35="

$tbl.Rows.Item(11).Cells.Item(1).Range.Text = "61+32="
$tbl.Rows.Item(11).Cells.Item(2).Range.Text = "70-30="
$tbl.Rows.Item(11).Cells.Item(3).Range.Text = "34+38="
$tbl.Rows.Item(11).Cells.Item(4).Range.Text = "37+5="
$tbl.Rows.Item(11).Cells.Item(5).Range.Text = "33-5="

$tbl.Rows.Item(12).Cells.Item(1).Range.Text = "35+23="
$tbl.Rows.Item(12).Cells.Item(2).Range.Text = "24-15="
$tbl.Rows.Item(12).Cells.Item(3).Range.Text = "2+0="
$tbl.Rows.Item(12).Cells.Item(4).Range.Text = "46-20="
$tbl.Rows.Item(12).Cells.Item(5).Range.Text = "78-61="

$tbl.Rows.Item(13).Cells.Item(1).Range.Text = "36+31="
$tbl.Rows.Item(13).Cells.Item(2).Range.Text = "94-77="
$tbl.Rows.Item(13).Cells.Item(3).Range.Text = "30-18="
$tbl.Rows.Item(13).Cells.Item(4).Range.Text = "47+25="
$tbl.Rows.Item(13).Cells.Item(5).Range.Text = "8-4="

$tbl.Rows.Item(14).Cells.Item(1).Range.Text = "88-81="
$tbl.Rows.Item(14).Cells.Item(2).Range.Text = "70-34="
$tbl.Rows.Item(14).Cells.Item(3).Range.Text = "52-4="
$tbl.Rows.Item(14).Cells.Item(4).Range.Text = "85-54="
$tbl.Rows.Item(14).Cells.Item(5).Range.Text = "54+26="

$tbl.Rows.Item(15).Cells.Item(1).Range.Text = "63+27="
$tbl.Rows.Item(15).Cells.Item(2).Range.Text = "58-24="
$tbl.Rows.Item(15).Cells.Item(3).Range.Text = "87-20="
$tbl.Rows.Item(15).Cells.Item(4).Range.Text = "41-20="
$tbl.Rows.Item(15).Cells.Item(5).Range.Text = "62+11="

$tbl.Rows.Item(16).Cells.Item(1).Range.Text = "43+50="
$tbl.Rows.Item(16).Cells.Item(2).Range.Text = "65-56="
$tbl.Rows.Item(16).Cells.Item(3).Range.Text = "50+48="
$tbl.Rows.Item(16).Cells.Item(4).Range.Text = "97-81="
$tbl.Rows.Item(16).Cells.Item(5).Range.Text = "14-7="

$tbl.Rows.Item(17).Cells.Item(1).Range.Text = "61-2="
$tbl.Rows.Item(17).Cells.Item(2).Range.Text = "96-12="
$tbl.Rows.Item(17).Cells.Item(3).Range.Text = "82-27="
$tbl.Rows.Item(17).Cells.Item(4).Range.Text = "11+72="
$tbl.Rows.Item(17).Cells.Item(5).Range.Text = "70-39="

$tbl.Rows.Item(18).Cells.Item(1).Range.Text = "36-24="
$tbl.Rows.Item(18).Cells.Item(2).Range.Text = "99-31="
$tbl.Rows.Item(18).Cells.Item(3).Range.Text = "63-13="
$tbl.Rows.Item(18).Cells.Item(4).Range.Text = "32+4="
$tbl.Rows.Item(18).Cells.Item(5).Range.Text = "96-86="

$tbl.Rows.Item(19).Cells.Item(1).Range.Text = "76-75="
$tbl.Rows.Item(19).Cells.Item(2).Range.Text = "83-48="
$tbl.Rows.Item(19).Cells.Item(3).Range.Text = "20+29="
$tbl.Rows.Item(19).Cells.Item(4).Range.Text = "29-5="
$tbl.Rows.Item(19).Cells.Item(5).Range.Text = "26+58="

$tbl.Rows.Item(20).Cells.Item(1).Range.Text = "75-69="
$tbl.Rows.Item(20).Cells.Item(2).Range.Text = "49-38="
$tbl.Rows.Item(20).Cells.Item(3).Range.Text = "89-87="
$tbl.Rows.Item(20).Cells.Item(4).Range.Text = "12+0="
$tbl.Rows.Item(20).Cells.Item(5).Range.Text = "91-81="
